$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before existing row 18 (PPC_P_SET_GRIDOP_ABS) for the new
# "PPC_P_SET_ABS" datapoint.
$ws.Rows.Item(18).Insert()
$ws.Range("A18").Value = "datapoints"
$ws.Range("B18").Value = "PPC_P_SET_ABS"
$ws.Range("C18").Value = "W"
$ws.Range("D18").Value = "Absolute active power setpoint"

# After the insert, the old row 19 (PPC_P_SET_GRIDOP_REL) becomes row 20.
# Insert two new rows after it (before the old PPC_P_SET_REL row, now at 21)
# for the "manually" absolute/relative setpoints.
$ws.Rows.Item(21).Insert()
$ws.Range("A21").Value = "datapoints"
$ws.Range("B21").Value = "PPC_P_SET_MANUAL_ABS"
$ws.Range("C21").Value = "W"
$ws.Range("D21").Value = "Absolute active power setpoint (manually)"

$ws.Rows.Item(22).Insert()
$ws.Range("A22").Value = "datapoints"
$ws.Range("B22").Value = "PPC_P_SET_MANUAL_REL"
$ws.Range("C22").Value = "%"
$ws.Range("D22").Value = "Relative active power setpoint (manually)"

# PPC_P_SET_REL is now row 23. Insert a new row after it (before
# PPC_P_SET_RPC_REL, now at 24) for the 3rd-party absolute setpoint.
$ws.Rows.Item(24).Insert()
$ws.Range("A24").Value = "datapoints"
$ws.Range("B24").Value = "PPC_P_SET_RPC_ABS"
$ws.Range("C24").Value = "W"
$ws.Range("D24").Value = "Absolute active power setpoint (3rd party)"
